$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 21154.2
$ws.Range("I111").Value = 1005.5
$ws.Range("K111").Value = 3016.5
$ws.Range("M111").Value = 50.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1523310.1
$ws.Range("I132").Value = 967.5074499999999
$ws.Range("J132").Value = 18522802
$ws.Range("K132").Value = 2902.52235
$ws.Range("L132").Value = 55568406
$ws.Range("M132").Value = -372.5223499999997
$ws.Range("N132").Value = -55573466

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 1000
$ws.Range("J3").Value = 1000
$ws.Range("L3").Value = 1000
$ws.Range("N3").Value = -1230

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 300
$ws.Range("J11").Value = 300
$ws.Range("L11").Value = 300
$ws.Range("N11").Value = -588

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 2994
$ws.Range("I14").Value = 2994
$ws.Range("K14").Value = 2994
$ws.Range("M14").Value = -2819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 869198
$ws.Range("I61").Value = 969947.6
$ws.Range("J61").Value = 2751.6
$ws.Range("K61").Value = 969947.6
$ws.Range("L61").Value = 2751.6
$ws.Range("M61").Value = -969735.6
$ws.Range("N61").Value = -3175.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 63890830
$ws.Range("I74").Value = 58824376
$ws.Range("J74").Value = 76195080
$ws.Range("K74").Value = 58824376
$ws.Range("L74").Value = 76195080
$ws.Range("M74").Value = -58823502
$ws.Range("N74").Value = -76196828

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 63890830
$ws.Range("I77").Value = 58824376
$ws.Range("J77").Value = 76195080
$ws.Range("K77").Value = 294121880
$ws.Range("L77").Value = 380975400
$ws.Range("M77").Value = -294117512
$ws.Range("N77").Value = -380984136

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 748.86365
$ws.Range("I110").Value = 770.6667
$ws.Range("J110").Value = 650.75
$ws.Range("K110").Value = 770.6667
$ws.Range("L110").Value = 650.75
$ws.Range("M110").Value = 1274.3333
$ws.Range("N110").Value = -4740.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1038
$ws.Range("I122").Value = 1201.9
$ws.Range("J122").Value = 491.66666
$ws.Range("K122").Value = 3605.7
$ws.Range("L122").Value = 1474.99998
$ws.Range("M122").Value = -1155.7
$ws.Range("N122").Value = -6374.999980000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 8683127
$ws.Range("I132").Value = 9618007
$ws.Range("J132").Value = 4631979
$ws.Range("K132").Value = 28854021
$ws.Range("L132").Value = 13895937
$ws.Range("M132").Value = -28851491
$ws.Range("N132").Value = -13900997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 869198
$ws.Range("I136").Value = 969947.6
$ws.Range("J136").Value = 2751.6
$ws.Range("K136").Value = 2909842.8
$ws.Range("L136").Value = 8254.799999999999
$ws.Range("M136").Value = -2907292.8
$ws.Range("N136").Value = -13354.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 24351608
$ws.Range("I134").Value = 23810528
$ws.Range("J134").Value = 35714284
$ws.Range("K134").Value = 71431584
$ws.Range("L134").Value = 107142852
$ws.Range("M134").Value = -71429049
$ws.Range("N134").Value = -107147922

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1528187.5
$ws.Range("I31").Value = 1848.0358
$ws.Range("K31").Value = 1848.0358
$ws.Range("M31").Value = -1553.0358

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1528187.5
$ws.Range("I34").Value = 1848.0358
$ws.Range("K34").Value = 1848.0358
$ws.Range("M34").Value = -1646.0358

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1467141.9
$ws.Range("I58").Value = 724.381
$ws.Range("K58").Value = 724.381
$ws.Range("M58").Value = -521.381

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3691.7693
$ws.Range("I122").Value = 4731.731
$ws.Range("J122").Value = 1611.8462
$ws.Range("K122").Value = 14195.193
$ws.Range("L122").Value = 4835.5386
$ws.Range("M122").Value = -11745.193
$ws.Range("N122").Value = -9735.5386

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1467141.9
$ws.Range("I136").Value = 724.381
$ws.Range("K136").Value = 2173.143
$ws.Range("M136").Value = 376.857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 157866.56
$ws.Range("J141").Value = 157866.56
$ws.Range("L141").Value = 157866.56
$ws.Range("N141").Value = -168226.56

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6486897
$ws.Range("I5").Value = 25642324
$ws.Range("J5").Value = 2382163
$ws.Range("K5").Value = 76926972
$ws.Range("L5").Value = 7146489
$ws.Range("M5").Value = -76926860
$ws.Range("N5").Value = -7146713

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 249.94118
$ws.Range("I26").Value = 59.85
$ws.Range("J26").Value = 521.5
$ws.Range("K26").Value = 179.55
$ws.Range("L26").Value = 1564.5
$ws.Range("M26").Value = 108.45
$ws.Range("N26").Value = -2140.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1893.9697
$ws.Range("I113").Value = 1477.5
$ws.Range("J113").Value = 2131.9524
$ws.Range("K113").Value = 4432.5
$ws.Range("L113").Value = 6395.8572
$ws.Range("M113").Value = -2262.5
$ws.Range("N113").Value = -10735.8572

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2327.862
$ws.Range("I132").Value = 2108.5715
$ws.Range("J132").Value = 2532.5334
$ws.Range("K132").Value = 18977.1435
$ws.Range("L132").Value = 22792.8006
$ws.Range("M132").Value = -16447.1435
$ws.Range("N132").Value = -27852.8006

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 6486897
$ws.Range("I135").Value = 25642324
$ws.Range("J135").Value = 2382163
$ws.Range("K135").Value = 230780916
$ws.Range("L135").Value = 21439467
$ws.Range("M135").Value = -230778381
$ws.Range("N135").Value = -21444537

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 20307.691
$ws.Range("I113").Value = 780
$ws.Range("K113").Value = 780
$ws.Range("M113").Value = 1390

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1926.2727
$ws.Range("I61").Value = 1864.8334
$ws.Range("K61").Value = 1864.8334
$ws.Range("M61").Value = -1662.8334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1926.2727
$ws.Range("I113").Value = 1864.8334
$ws.Range("K113").Value = 1864.8334
$ws.Range("M113").Value = 305.1666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 9623039
$ws.Range("I122").Value = 1332187.5
$ws.Range("J122").Value = 28573556
$ws.Range("K122").Value = 3996562.5
$ws.Range("L122").Value = 85720668
$ws.Range("M122").Value = -3994112.5
$ws.Range("N122").Value = -85725568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3488327.2
$ws.Range("I132").Value = 4468339
$ws.Range("J132").Value = 3840.3333
$ws.Range("K132").Value = 13405017
$ws.Range("L132").Value = 11520.9999
$ws.Range("M132").Value = -13402487
$ws.Range("N132").Value = -16580.9999
